$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row just above the current row 219 (a new week of
# "Zapallo italiano" price data for Región de O'Higgins). This shifts the
# previous rows 219..270 down to 220..271, matching the target diff.
$ws.Rows(219).Insert()

$ws.Range("A219").Value = 9
$ws.Range("B219").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C219").Value = "Metropolitana"
$ws.Range("D219").Value = 44551
$ws.Range("E219").Value = 13
$ws.Range("F219").Value = 100112032
$ws.Range("G219").Value = "Zapallo italiano"
$ws.Range("H219").Value = "Sin especificar"
$ws.Range("I219").Value = "Primera"
$ws.Range("J219").Value = 106
$ws.Range("K219").Value = 8000
$ws.Range("L219").Value = 9000
$ws.Range("M219").Value = 8500
$ws.Range("N219").Value = "`$/caja 50 unidades"
$ws.Range("O219").Value = "Región de O'Higgins"
$ws.Range("P219").Value = 170
$ws.Range("Q219").Value = 50
$ws.Range("R219").Value = "Hortaliza"
